$wb = $excel.ActiveWorkbook

# --- Flight Mission Cycle: update cycle counts / setting labels ---
$wsFMC = $wb.Worksheets.Item("Flight Mission Cycle")
$wsFMC.Range("B2").Value = 2
$wsFMC.Range("A4").Value = "Piano"

# --- Typing: Force is now angle-dependent (single set point) instead of a set_points curve ---
$wsTyping = $wb.Worksheets.Item("Typing")
$wsTyping.Range("B2").Value = "angle_dependent"
$wsTyping.Range("C2").Value = "<1"
$wsTyping.Range("D2").Value = 20
$wsTyping.Range("E2").ClearContents()
$wsTyping.Range("C3").Value = 30
$wsTyping.Range("D3:E3").ClearContents()

# --- Light switch: Force is now angle-dependent (single set point) instead of a set_points curve ---
$wsLS = $wb.Worksheets.Item("Light switch")
$wsLS.Range("B2").Value = "angle_dependent"
$wsLS.Range("C2").Value = ">10"
$wsLS.Range("D2").Value = 30
$wsLS.Range("E2:F2").ClearContents()
$wsLS.Range("C3").Value = 40
$wsLS.Range("D3:F3").ClearContents()

# --- Piano: bring layout in line with the other settings sheets (Type column + RoM row) ---
$wsPiano = $wb.Worksheets.Item("Piano")
$wsPiano.Range("B1:B3").Insert(-4161)
$wsPiano.Range("A4:H4").Insert(-4121)
$wsPiano.Range("B1").Value = "Type"
$wsPiano.Range("B2").Value = "set_points"
$wsPiano.Range("A4").Value = "RoM"
$wsPiano.Range("B4").Value = "triangle"

# --- Restore/record the cursor position on each touched sheet ---
$wsFMC.Range("H9").Select()
$wsTyping.Range("F7").Select()
$wsPiano.Range("I10").Select()
$wsLS.Range("E10").Select()
